$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended after the existing last row (173).
# Columns: A = Serie (date label), C = "3 en 2", D = "5 en 5"
$ws.Range("A174").Value = "27-09-2021"
$ws.Range("D174").Value = 3.16

$ws.Range("A175").Value = "28-09-2021"
$ws.Range("C175").Value = 2.86
$ws.Range("D175").Value = 3.21

$ws.Range("A176").Value = "29-09-2021"
$ws.Range("C176").Value = 3.03
$ws.Range("D176").Value = 3.28

$ws.Range("A177").Value = "30-09-2021"
$ws.Range("C177").Value = 3.16
$ws.Range("D177").Value = 3.26
